$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1470.48
$ws.Range("I40").Value = 1397.0769
$ws.Range("J40").Value = 1550
$ws.Range("K40").Value = 1397.0769
$ws.Range("L40").Value = 1550
$ws.Range("M40").Value = -1222.0769
$ws.Range("N40").Value = -1900
$ws.Range("H137").Value = 1254.8113
$ws.Range("I137").Value = 1095.7894
$ws.Range("J137").Value = 1657.6666
$ws.Range("K137").Value = 3287.3682
$ws.Range("L137").Value = 4972.9998
$ws.Range("M137").Value = -737.3681999999999
$ws.Range("N137").Value = -10072.9998
$ws.Range("H138").Value = 2921.8684
$ws.Range("J138").Value = 4484.793
$ws.Range("L138").Value = 13454.379
$ws.Range("N138").Value = -23734.379
$ws.Range("H141").Value = 5382.294
$ws.Range("I141").Value = 2496.5173
$ws.Range("J141").Value = 22119.8
$ws.Range("K141").Value = 7489.5519
$ws.Range("L141").Value = 66359.39999999999
$ws.Range("M141").Value = -2309.5519
$ws.Range("N141").Value = -76719.39999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("H61").Value = 830.44183
$ws.Range("I61").Value = 736.05884
$ws.Range("J61").Value = 1187
$ws.Range("K61").Value = 736.05884
$ws.Range("L61").Value = 1187
$ws.Range("M61").Value = -524.05884
$ws.Range("N61").Value = -1611
$ws.Range("H74").Value = 761.2558
$ws.Range("I74").Value = 718.7027
$ws.Range("K74").Value = 718.7027
$ws.Range("M74").Value = 155.2973
$ws.Range("H77").Value = 761.2558
$ws.Range("I77").Value = 718.7027
$ws.Range("K77").Value = 3593.5135
$ws.Range("M77").Value = 774.4865
$ws.Range("H124").Value = 25119.334
$ws.Range("J124").Value = 25119.334
$ws.Range("L124").Value = 25119.334
$ws.Range("N124").Value = -34939.334
$ws.Range("H125").Value = 70713.57000000001
$ws.Range("J125").Value = 70713.57000000001
$ws.Range("L125").Value = 70713.57000000001
$ws.Range("N125").Value = -80553.57000000001
$ws.Range("H132").Value = 1297.6786
$ws.Range("I132").Value = 850.0833
$ws.Range("J132").Value = 2103.35
$ws.Range("K132").Value = 2550.2499
$ws.Range("L132").Value = 6310.049999999999
$ws.Range("M132").Value = -20.2498999999998
$ws.Range("N132").Value = -11370.05
$ws.Range("H136").Value = 830.44183
$ws.Range("I136").Value = 736.05884
$ws.Range("J136").Value = 1187
$ws.Range("K136").Value = 2208.17652
$ws.Range("L136").Value = 3561
$ws.Range("M136").Value = 341.82348
$ws.Range("N136").Value = -8661
$ws.Range("M5").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("H22").Value = 8966.333000000001
$ws.Range("I22").Value = 10543.6
$ws.Range("J22").Value = 1080
$ws.Range("K22").Value = 10543.6
$ws.Range("L22").Value = 1080
$ws.Range("M22").Value = -10370.6
$ws.Range("N22").Value = -1426
$ws.Range("H134").Value = 1675.0754
$ws.Range("I134").Value = 1511.1395
$ws.Range("J134").Value = 2380
$ws.Range("K134").Value = 4533.4185
$ws.Range("L134").Value = 7140
$ws.Range("M134").Value = -1998.4185
$ws.Range("N134").Value = -12210
$ws.Range("M4").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1808.2322
$ws.Range("I31").Value = 1302.921
$ws.Range("J31").Value = 2875
$ws.Range("K31").Value = 1302.921
$ws.Range("L31").Value = 2875
$ws.Range("M31").Value = -1007.921
$ws.Range("N31").Value = -3465
$ws.Range("H34").Value = 1808.2322
$ws.Range("I34").Value = 1302.921
$ws.Range("J34").Value = 2875
$ws.Range("K34").Value = 1302.921
$ws.Range("L34").Value = 2875
$ws.Range("M34").Value = -1100.921
$ws.Range("N34").Value = -3279
$ws.Range("H58").Value = 1002388.44
$ws.Range("I58").Value = 1544142
$ws.Range("K58").Value = 1544142
$ws.Range("M58").Value = -1543939
$ws.Range("H132").Value = 226178.83
$ws.Range("I132").Value = 265726.6
$ws.Range("J132").Value = 2074.889
$ws.Range("K132").Value = 797179.7999999999
$ws.Range("L132").Value = 6224.667
$ws.Range("M132").Value = -794649.7999999999
$ws.Range("N132").Value = -11284.667
$ws.Range("H134").Value = 1049.6624
$ws.Range("I134").Value = 823.4737
$ws.Range("J134").Value = 1694.3
$ws.Range("K134").Value = 2470.4211
$ws.Range("L134").Value = 5082.9
$ws.Range("M134").Value = 64.57889999999998
$ws.Range("N134").Value = -10152.9
$ws.Range("H136").Value = 1002388.44
$ws.Range("I136").Value = 1544142
$ws.Range("K136").Value = 4632426
$ws.Range("M136").Value = -4629876

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1175.697
$ws.Range("I5").Value = 1185.3334
$ws.Range("J5").Value = 1150
$ws.Range("K5").Value = 3556.0002
$ws.Range("L5").Value = 3450
$ws.Range("M5").Value = -3444.0002
$ws.Range("N5").Value = -3674
$ws.Range("H12").Value = 690176.9
$ws.Range("I12").Value = 80.85714
$ws.Range("J12").Value = 920208.9
$ws.Range("K12").Value = 242.57142
$ws.Range("L12").Value = 2760626.7
$ws.Range("M12").Value = -69.57141999999999
$ws.Range("N12").Value = -2760972.7
$ws.Range("H122").Value = 681.0417
$ws.Range("I122").Value = 498.75
$ws.Range("J122").Value = 863.3333
$ws.Range("K122").Value = 4488.75
$ws.Range("L122").Value = 7769.9997
$ws.Range("M122").Value = -2038.75
$ws.Range("N122").Value = -12669.9997
$ws.Range("H135").Value = 1175.697
$ws.Range("I135").Value = 1185.3334
$ws.Range("J135").Value = 1150
$ws.Range("K135").Value = 10668.0006
$ws.Range("L135").Value = 10350
$ws.Range("M135").Value = -8133.000599999999
$ws.Range("N135").Value = -15420

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2077.1304
$ws.Range("I132").Value = 1119.75
$ws.Range("J132").Value = 3121.5454
$ws.Range("K132").Value = 3359.25
$ws.Range("L132").Value = 9364.636200000001
$ws.Range("M132").Value = -829.25
$ws.Range("N132").Value = -14424.6362

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6086.385
$ws.Range("I7").Value = 6667.6665
$ws.Range("J7").Value = 5912
$ws.Range("K7").Value = 6667.6665
$ws.Range("L7").Value = 5912
$ws.Range("M7").Value = -6555.6665
$ws.Range("N7").Value = -6136
$ws.Range("H16").Value = 2907.9
$ws.Range("I16").Value = 3064.3333
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 3064.3333
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -2894.3333
$ws.Range("N16").Value = -1840
$ws.Range("H22").Value = 1120
$ws.Range("I22").Value = 760
$ws.Range("J22").Value = 1248.5714
$ws.Range("K22").Value = 760
$ws.Range("L22").Value = 1248.5714
$ws.Range("M22").Value = -465
$ws.Range("N22").Value = -1838.5714
$ws.Range("H27").Value = 1120
$ws.Range("I27").Value = 760
$ws.Range("J27").Value = 1248.5714
$ws.Range("K27").Value = 760
$ws.Range("L27").Value = 1248.5714
$ws.Range("M27").Value = -653
$ws.Range("N27").Value = -1462.5714
$ws.Range("H43").Value = 27500
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("N43").Value = -5386
$ws.Range("H46").Value = 1608.1666
$ws.Range("I46").Value = 1413.7142
$ws.Range("J46").Value = 1880.4
$ws.Range("K46").Value = 1413.7142
$ws.Range("L46").Value = 1880.4
$ws.Range("M46").Value = -1225.7142
$ws.Range("N46").Value = -2256.4
$ws.Range("H126").Value = 6086.385
$ws.Range("I126").Value = 6667.6665
$ws.Range("J126").Value = 5912
$ws.Range("K126").Value = 20002.9995
$ws.Range("L126").Value = 17736
$ws.Range("M126").Value = -17532.9995
$ws.Range("N126").Value = -22676
$ws.Range("H132").Value = 2687.7637
$ws.Range("I132").Value = 2096.0652
$ws.Range("K132").Value = 6288.1956
$ws.Range("M132").Value = -3758.1956
$ws.Range("H136").Value = 1349.9524
$ws.Range("I136").Value = 1084.7354
$ws.Range("K136").Value = 3254.2062
$ws.Range("M136").Value = -704.2062000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 761.0789
$ws.Range("I132").Value = 690.7778
$ws.Range("J132").Value = 933.63635
$ws.Range("K132").Value = 2072.3334
$ws.Range("L132").Value = 2800.90905
$ws.Range("M132").Value = 457.6666
$ws.Range("N132").Value = -7860.90905
$ws.Range("H136").Value = 787.53125
$ws.Range("I136").Value = 784.38464
$ws.Range("K136").Value = 2353.15392
$ws.Range("M136").Value = 196.8460800000003

Write-Host "Applied all updates"
